# Apply quarterly-style price/volume refresh to the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several "Price" values look numeric (e.g. "141.50") but must stay as literal
# text so trailing/insignificant-looking digits are not dropped by Excel's
# automatic number detection. Force those specific cells to Text format first.
$textPriceCells = @("D5", "D6", "D7", "D9", "D10", "D14", "D17", "D20", "D21", "D22", "D23", "D25", "D28", "D32", "D37", "D40", "D42", "D44", "D45", "D47", "D48", "D49", "D51")
foreach ($c in $textPriceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.191.26'
$ws.Range("E2").Value = '  +3.71%  '
$ws.Range("D3").Value = '1.603.06'
$ws.Range("E3").Value = '  +2.67%  '
$ws.Range("E4").Value = '  -0.43%  '
$ws.Range("D5").Value = '212.78'
$ws.Range("E5").Value = '  +2.88%  '
$ws.Range("D6").Value = '0.999'
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("D7").Value = '0.485'
$ws.Range("E7").Value = '  +1.61%  '
$ws.Range("E8").Value = '  +3.36%  '
$ws.Range("D9").Value = '0.0617'
$ws.Range("E9").Value = '  +1.94%  '
$ws.Range("D10").Value = '18.02'
$ws.Range("E10").Value = '  +1.67%  '
$ws.Range("E11").Value = '  +4.58%  '
$ws.Range("D13").Value = '1.602.14'
$ws.Range("E13").Value = '  +2.58%  '
$ws.Range("D14").Value = '4.01'
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("E15").Value = '  +1.67%  '
$ws.Range("D16").Value = '26.104.48'
$ws.Range("E16").Value = '  +3.32%  '
$ws.Range("D17").Value = '60.51'
$ws.Range("E17").Value = '  +2.35%  '
$ws.Range("D18").Value = '0.0₃0722'
$ws.Range("E18").Value = '  +1.78%  '
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = '204.85'
$ws.Range("E20").Value = '  +10.60%  '
$ws.Range("D21").Value = '4.25'
$ws.Range("E21").Value = '  +3.43%  '
$ws.Range("D22").Value = '9.33'
$ws.Range("E22").Value = '  +0.86%  '
$ws.Range("D23").Value = '5.99'
$ws.Range("E23").Value = '  +2.16%  '
$ws.Range("E24").Value = '  +11.17%  '
$ws.Range("D25").Value = '141.50'
$ws.Range("E25").Value = '  +1.38%  '
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("E27").Value = '  -2.40%  '
$ws.Range("D28").Value = '15.24'
$ws.Range("E28").Value = '  +2.99%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  +1.68%  '
$ws.Range("E31").Value = '  +1.77%  '
$ws.Range("D32").Value = '3.13'
$ws.Range("E32").Value = '  +3.50%  '
$ws.Range("E33").Value = '  -0.52%  '
$ws.Range("E34").Value = '  +1.86%  '
$ws.Range("E35").Value = '  +2.00%  '
$ws.Range("D36").Value = '1.112.45'
$ws.Range("E36").Value = '  +2.32%  '
$ws.Range("D37").Value = '0.0161'
$ws.Range("E37").Value = '  +7.88%  '
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("D40").Value = '0.779'
$ws.Range("E40").Value = '  +2.43%  '
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("D42").Value = '0.782'
$ws.Range("E42").Value = '  -5.25%  '
$ws.Range("D43").Value = '1.733.82'
$ws.Range("E43").Value = '  +2.33%  '
$ws.Range("D44").Value = '92.85'
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").Value = '5.11'
$ws.Range("E45").Value = '  +1.04%  '
$ws.Range("E46").Value = '  +5.33%  '
$ws.Range("D47").Value = '53.58'
$ws.Range("E47").Value = '  +2.31%  '
$ws.Range("D48").Value = '0.0503'
$ws.Range("E48").Value = '  -0.17%  '
$ws.Range("D49").Value = '0.409'
$ws.Range("E49").Value = '  +0.71%  '
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("D51").Value = '7.25'
$ws.Range("E51").Value = '  +1.93%  '
